# Apply the weekly-data commit: backfill Q/O/R fixes on existing rows and
# append 26 new weekly OHLCV rows (113-138) for LICI.NS, mirroring the
# "break out stock.yaml completed" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to already-present rows ---
$ws.Range("Q65").Value  = 0
$ws.Range("Q76").Value  = 0
$ws.Range("O109").Value = 2
$ws.Range("R111").Value = 0
$ws.Range("R112").Value = 0

# --- New weekly rows 113-138 ---
# Columns: A Datetime | B Open | C High | D Low | E Close | F Adj Close |
#          G Volume | H Year | I Month | J Day | K Hour | L Minute |
#          M Second | N Week | O isPivot | P two_line_structure |
#          Q detect_structure | R backup
# "__BLANK__" marks columns that stay empty (Adj Close / backup), matching
# the source rows which only carry an (empty) placeholder cell there.
$newRows = @(
    @(45474, 989.6199570803844, 1018.960448270205, 976.8891798237944, 1010.158325195312, "__BLANK__", 10704417, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0, "__BLANK__"),
    @(45481, 1011.501023066534, 1060.186299063986, 993.6977851041546, 1045.8642578125, "__BLANK__", 22786039, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0, "__BLANK__"),
    @(45488, 1049.295497062162, 1123.889963677955, 1032.984149413455, 1102.257568359375, "__BLANK__", 30213375, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1, "__BLANK__"),
    @(45495, 1100, 1197, 1044.5, 1184.349975585938, "__BLANK__", 34438020, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0, "__BLANK__"),
    @(45502, 1194.650024414062, 1222, 1160, 1179.550048828125, "__BLANK__", 19154948, 2024, 7, 29, 0, 0, 0, 31, 1, 0, 0, "__BLANK__"),
    @(45509, 1154, 1160, 1076.199951171875, 1133.449951171875, "__BLANK__", 17925128, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0, "__BLANK__"),
    @(45516, 1125, 1130.300048828125, 1003.75, 1057.300048828125, "__BLANK__", 15708123, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0, "__BLANK__"),
    @(45523, 1059.650024414062, 1092, 1056.650024414062, 1061.349975585938, "__BLANK__", 6248392, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0, "__BLANK__"),
    @(45530, 1065, 1095, 1045, 1064.25, "__BLANK__", 8208283, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0, "__BLANK__"),
    @(45537, 1065.5, 1079.949951171875, 1035, 1038.900024414062, "__BLANK__", 5083514, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0, "__BLANK__"),
    @(45544, 1038.900024414062, 1046, 1004.849975585938, 1034.150024414062, "__BLANK__", 4865139, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 2, "__BLANK__"),
    @(45551, 1039, 1041.650024414062, 991, 1010.349975585938, "__BLANK__", 5068734, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0, "__BLANK__"),
    @(45558, 1012.950012207031, 1048.900024414062, 1010.5, 1024.25, "__BLANK__", 7096412, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0, "__BLANK__"),
    @(45565, 1022, 1022, 955.7000122070312, 971.1500244140625, "__BLANK__", 5889264, 2024, 9, 30, 0, 0, 0, 40, 0, 0, 0, "__BLANK__"),
    @(45572, 971, 979.2000122070312, 923.0999755859375, 948.8499755859375, "__BLANK__", 8844125, 2024, 10, 7, 0, 0, 0, 41, 0, 0, 0, "__BLANK__"),
    @(45579, 958.2999877929688, 965.9500122070312, 925.2999877929688, 934.9000244140625, "__BLANK__", 4273247, 2024, 10, 14, 0, 0, 0, 42, 0, 0, 0, "__BLANK__"),
    @(45586, 945, 952.9000244140625, 888.2999877929688, 903.6500244140625, "__BLANK__", 5398572, 2024, 10, 21, 0, 0, 0, 43, 0, 0, 0, "__BLANK__"),
    @(45593, 909.3499755859375, 947.4500122070312, 896.6500244140625, 930.3499755859375, "__BLANK__", 3585636, 2024, 10, 28, 0, 0, 0, 44, 0, 0, 0, "__BLANK__"),
    @(45600, 934, 954, 911.25, 914.6500244140625, "__BLANK__", 4010310, 2024, 11, 4, 0, 0, 0, 45, 0, 0, 0, "__BLANK__"),
    @(45607, 920, 958, 890.0499877929688, 908.6500244140625, "__BLANK__", 6042877, 2024, 11, 11, 0, 0, 0, 46, 0, 0, 0, "__BLANK__"),
    @(45614, 908.6500244140625, 916.9500122070312, 872, 889.5, "__BLANK__", 4041182, 2024, 11, 18, 0, 0, 0, 47, 2, 0, 0, "__BLANK__"),
    @(45621, 905.9500122070312, 992.1500244140625, 900, 985.5, "__BLANK__", 12956216, 2024, 11, 25, 0, 0, 0, 48, 0, 0, 0, "__BLANK__"),
    @(45628, 985, 999.9500122070312, 966.2000122070312, 983.25, "__BLANK__", 7224787, 2024, 12, 2, 0, 0, 0, 49, 0, 0, 0, "__BLANK__"),
    @(45635, 987.6500244140625, 1007.799987792969, 913.5499877929688, 931.9000244140625, "__BLANK__", 9436636, 2024, 12, 9, 0, 0, 0, 50, 0, 0, 0, "__BLANK__"),
    @(45642, 935, 936.5, 890.5499877929688, 901.7000122070312, "__BLANK__", 5764231, 2024, 12, 16, 0, 0, 0, 51, 0, 0, 0, "__BLANK__"),
    @(45649, 902.5, 906.9000244140625, 885.1500244140625, 887, "__BLANK__", 4286778, 2024, 12, 23, 0, 0, 0, 52, 0, 0, 0, "__BLANK__")
)

$startRow = 113
$r = $startRow
foreach ($row in $newRows) {
    $c = 1
    foreach ($v in $row) {
        if ($v -ne "__BLANK__") {
            $cell = $ws.Cells.Item($r, $c)
            $cell.Value = $v
            if ($c -eq 1) {
                $cell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
            }
        }
        $c = $c + 1
    }
    $r = $r + 1
}
